$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data table to write: File Number | University | Group Type
$table = @(
    @("File Number", "University", "Group Type"),
    @(1,  "СПбГУТ",  "Основная"),
    @(2,  "СПбГУТ",  "Подслушано"),
    @(3,  "СПбГУТ",  "Поступление"),
    @(4,  "НИУ ВШЭ", "Основная"),
    @(5,  "НИУ ВШЭ", "Подслушано"),
    @(6,  "НИУ ВШЭ", "Поступление"),
    @(7,  "ИТМО",    "Основная"),
    @(8,  "ИТМО",    "Подслушано"),
    @(9,  "ИТМО",    "Поступление"),
    @(10, "ЛЭТИ",    "Основная"),
    @(11, "ЛЭТИ",    "Подслушано"),
    @(12, "ЛЭТИ",    "Поступление"),
    @(13, "СПбПУ",   "Основная"),
    @(14, "СПбПУ",   "Подслушано"),
    @(15, "СПбПУ",   "Поступление")
)

# Remove the rows that are no longer part of the table (old data went to row 29)
$oldLastRow = 29
$newLastRow = $table.Count
if ($oldLastRow -gt $newLastRow) {
    $deleteRange = $ws.Range("A$($newLastRow + 1):A$oldLastRow")
    $deleteRange.EntireRow.Delete() | Out-Null
}

# Write the new table contents
for ($i = 0; $i -lt $table.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $table[$i][0]
    $ws.Cells.Item($r, 2).Value = $table[$i][1]
    $ws.Cells.Item($r, 3).Value = $table[$i][2]
}

# Update the selected cell to match the saved view
$ws.Range("E16").Select() | Out-Null
